$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new row (row 54) continuing the existing forecast series table.
$row = 54

# Column A: date value, formatted like the preceding date cells (style copied from A53)
$ws.Cells.Item($row - 1, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)
$ws.Cells.Item($row, 1).Value = 45986

# Column B: starting year of the forecast window
$ws.Cells.Item($row, 2).Value = 2025

# Column C: year-over-year forecast value
$ws.Cells.Item($row, 3).Value = 2.560577522109297

# Column D: ending year of the forecast window
$ws.Cells.Item($row, 4).Value = 2026

# Column E: year-over-year forecast value
$ws.Cells.Item($row, 5).Value = 2.991302072731838

$wb.Save()
